# Update stats for 2026-01 (row 26 of Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B26").Value = 6532
$ws.Range("D26").Value = 6087506
$ws.Range("E26").Value = 931.9513165952235
$ws.Range("F26").Value = 10.39378063207708
$ws.Range("H26").Value = 26.77251753624024
